$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.973.31'
$ws.Range("E2").Value = '  +1.18%  '
$ws.Range("D3").Value = '3.751.49'
$ws.Range("E3").Value = '  +2.17%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.15'
$ws.Range("E5").Value = '  +0.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.38'
$ws.Range("E6").Value = '  +0.88%  '
$ws.Range("D7").Value = '3.749.51'
$ws.Range("E7").Value = '  +2.24%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  +1.35%  '
$ws.Range("E10").Value = '  +2.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.48'
$ws.Range("E11").Value = '  +3.10%  '
$ws.Range("E12").Value = '  +0.72%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.86'
$ws.Range("E13").Value = '  +0.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000249'
$ws.Range("E14").Value = '  +1.85%  '
$ws.Range("D15").Value = '4.383.01'
$ws.Range("E15").Value = '  +2.26%  '
$ws.Range("D16").Value = '3.752.56'
$ws.Range("D17").Value = '69.037.53'
$ws.Range("E17").Value = '  +1.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.36'
$ws.Range("E18").Value = '  +2.08%  '
$ws.Range("E19").Value = '  -0.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.24'
$ws.Range("E20").Value = '  +1.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.83'
$ws.Range("E21").Value = '  +19.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '492.79'
$ws.Range("E22").Value = '  +0.46%  '
$ws.Range("E23").Value = '  +0.97%  '
$ws.Range("E24").Value = '  +8.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.79'
$ws.Range("E26").Value = '  +1.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.34'
$ws.Range("E27").Value = '  +1.56%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.12'
$ws.Range("E28").Value = '  +1.04%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("E30").Value = '  +3.37%  '
$ws.Range("E31").Value = '  +4.43%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.01'
$ws.Range("E32").Value = '  +2.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.55'
$ws.Range("E33").Value = '  +1.07%  '
$ws.Range("D34").Value = '3.897.53'
$ws.Range("E34").Value = '  +2.25%  '
$ws.Range("E35").Value = '  +0.44%  '
$ws.Range("D36").Value = '3.686.71'
$ws.Range("E36").Value = '  +1.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("E38").Value = '  +3.04%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.90'
$ws.Range("E39").Value = '  +3.14%  '
$ws.Range("E40").Value = '  +2.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.324'
$ws.Range("E41").Value = '  +1.47%  '
$ws.Range("E42").Value = '  +5.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '429.02'
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("E44").Value = '  +2.96%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '48.61'
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("E46").Value = '  +1.73%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.38'
$ws.Range("E48").Value = '  +0.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.55'
$ws.Range("E49").Value = '  +0.31%  '
$ws.Range("D50").Value = '2.789.47'
$ws.Range("E50").Value = '  +2.43%  '
$ws.Range("E51").Value = '  +1.27%  '
